$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Text = "28÷4=7, 0"
$cell = $t.Cell(1, 2)
$cell.Range.Text = "85÷8=10, 5"
$cell = $t.Cell(1, 3)
$cell.Range.Text = "11÷3=3, 2"
$cell = $t.Cell(1, 4)
$cell.Range.Text = "85÷9=9, 4"
$cell = $t.Cell(1, 5)
$cell.Range.Text = "53÷2=26, 1"
$cell = $t.Cell(5, 1)
$cell.Range.Text = "23÷4=5, 3"
$cell = $t.Cell(5, 2)
$cell.Range.Text = "21÷9=2, 3"
$cell = $t.Cell(5, 3)
$cell.Range.Text = "48÷8=6, 0"
$cell = $t.Cell(5, 4)
$cell.Range.Text = "56÷5=11, 1"
$cell = $t.Cell(5, 5)
$cell.Range.Text = "23÷3=7, 2"
$cell = $t.Cell(9, 1)
$cell.Range.Text = "91÷6=15, 1"
$cell = $t.Cell(9, 2)
$cell.Range.Text = "63÷5=12, 3"
$cell = $t.Cell(9, 3)
$cell.Range.Text = "49÷4=12, 1"
$cell = $t.Cell(9, 4)
$cell.Range.Text = "51÷3=17, 0"
$cell = $t.Cell(9, 5)
$cell.Range.Text = "72÷7=10, 2"
$cell = $t.Cell(13, 1)
$cell.Range.Text = "50÷9=5, 5"
$cell = $t.Cell(13, 2)
$cell.Range.Text = "99÷3=33, 0"
$cell = $t.Cell(13, 3)
$cell.Range.Text = "78÷6=13, 0"
$cell = $t.Cell(13, 4)
$cell.Range.Text = "47÷9=5, 2"
$cell = $t.Cell(13, 5)
$cell.Range.Text = "29÷8=3, 5"
$cell = $t.Cell(17, 1)
$cell.Range.Text = "35÷9=3, 8"
$cell = $t.Cell(17, 2)
$cell.Range.Text = "71÷6=11, 5"
$cell = $t.Cell(17, 3)
$cell.Range.Text = "14÷2=7, 0"
$cell = $t.Cell(17, 4)
$cell.Range.Text = "73÷5=14, 3"
$cell = $t.Cell(17, 5)
$cell.Range.Text = "15÷3=5, 0"

Write-Output "done"
